# Hale Aircraft Model, with bugs...
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Changchuan")

# ---------------------------------------------------------------------
# 1. New "Theta - Fig 12" results block added to the right of the
#    existing tables (columns K:N), mirroring the existing x/y,
#    Linear/Nonlinear layout used by the other blocks on this sheet.
# ---------------------------------------------------------------------
$ws1.Range("K1").Value = "Theta - Fig 12"
$ws1.Range("K1:N1").Merge()
$ws1.Range("K1:N1").HorizontalAlignment = -4108
$ws1.Range("K1:N1").VerticalAlignment = -4108

$ws1.Range("K2").Value = "Linear"
$ws1.Range("K2:L2").Merge()
$ws1.Range("M2").Value = "Nonlinear"
$ws1.Range("M2:N2").Merge()
$ws1.Range("K2:N2").HorizontalAlignment = -4108

$ws1.Range("K3").Value = "x"
$ws1.Range("L3").Value = "y"
$ws1.Range("M3").Value = "x"
$ws1.Range("N3").Value = "y"
$ws1.Range("K3:N3").HorizontalAlignment = -4108
$ws1.Range("K3:N3").VerticalAlignment = -4108

$thetaData = @(
    @(-0.00054059189732410997, -0.0097545206178084598, -0.00016849199663021501, -0.0077519379844979098),
    @(0.0284561410668633, 0.30835228755087302, 0.0309264692831769, 0.093610153562735804),
    @(0.091364400643624305, 0.911593486086968, 0.093391589795311897, 0.29511264097943701),
    @(0.15432565063061701, 1.5015108857839701, 0.15417873306632801, 0.495982040372211),
    @(0.216813859392778, 2.0708427344090001, 0.21587844423347199, 0.68265697339869902),
    @(0.27853492387909301, 2.6199874583823202, 0.278288219952952, 0.87809335843019198),
    @(0.34074496145238198, 3.1283238003995502, 0.34058910128090902, 1.0393011788575699),
    @(0.40264070638221899, 3.6052850703218202, 0.40243718004435802, 1.2091869762580001),
    @(0.46464186314136702, 4.0368116567906904, 0.46474046791703699, 1.3558471102584999),
    @(0.52679922104571397, 4.4232664406852003, 0.52655638452958198, 1.4797768064685499),
    @(0.58871449422021804, 4.7645017203981901, 0.58880624747128696, 1.5967062671608401),
    @(0.650571181033201, 5.0542251253362496, 0.65047668537047898, 1.6970377549125999),
    @(0.71234444921015505, 5.2951856850742596, 0.71210151660124299, 1.7887996495569101),
    @(0.77493384211308802, 5.4784550657226401, 0.77525586792704804, 1.83502500561267),
    @(0.83687765217314602, 5.6068918855736198, 0.83679201825480098, 1.88932528975078),
    @(0.89843918869689898, 5.6828408672951598, 0.89839767541163496, 1.9084388614677701),
    @(0.96075901852993295, 5.7124687643104597, 0.95989257168615505, 1.9290320334636799)
)

$r = 4
foreach ($row in $thetaData) {
    $ws1.Range("K$r").Value = $row[0]
    $ws1.Range("L$r").Value = $row[1]
    $ws1.Range("M$r").Value = $row[2]
    $ws1.Range("N$r").Value = $row[3]
    $r++
}

# A stray formatted (but emptied) cell left behind at the end of the
# pasted block, carrying the new underlined font.
$ws1.Range("L25").Font.Underline = $true

# ---------------------------------------------------------------------
# 2. Page setup for the data sheet.
# ---------------------------------------------------------------------
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# 3. New blank "Plan1" worksheet, added after "Changchuan".
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Plan1"
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1
$ws2.Range("A1:E27").Select()

# Restore the original sheet as the active/selected one and leave the
# cursor where the editor last left it.
$ws1.Activate()
$ws1.Range("F5").Select()
